# Insert a new weekly price record above the current row 71 ("Feria
# Lagunitas de Puerto Montt" - Poroto verde), pushing the existing
# rows 71-87 down to 72-88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 71:87 down to 72:88, creating a blank row 71.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new weekly record.
$ws.Range("A71").Value = 4
$ws.Range("B71").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C71").Value = "Los Lagos"
$ws.Range("D71").Value = 44754
$ws.Range("E71").Value = 10
$ws.Range("F71").Value = 100112031
$ws.Range("G71").Value = "Poroto verde"
$ws.Range("H71").Value = "Magnum"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 40
$ws.Range("K71").Value = 40000
$ws.Range("L71").Value = 40000
$ws.Range("M71").Value = 40000
$ws.Range("N71").Value = "$/malla 25 kilos"
$ws.Range("O71").Value = "Perú"
$ws.Range("P71").Value = 1600
$ws.Range("Q71").Value = 25
$ws.Range("R71").Value = "Hortaliza"
